# Auto-generated Excel COM-interop edit script
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Rename shared string "Bilan en Watt" -> "Bilan en kWh" (cell A13) ---
$ws.Range("A13").Value = "Bilan en kWh"

# --- 2) Update changed numeric cell values ---
$ws.Range("B2").Value = 43.252249734683943
$ws.Range("C2").Value = -25.437727650329101
$ws.Range("D2").Value = -2.4688041806060941
$ws.Range("BF2").Value = 0.35138171124437168
$ws.Range("BG2").Value = -27.526540443610511
$ws.Range("BH2").Value = -46.83088653208938
$ws.Range("B3").Value = -23.677387199202549
$ws.Range("C3").Value = -44.776091135500607
$ws.Range("D3").Value = -33.055579174106512
$ws.Range("I3").Value = -19.8814642654865
$ws.Range("J3").Value = -25.939295188562991
$ws.Range("K3").Value = -34.081066074403097
$ws.Range("S3").Value = -17.60574457490425
$ws.Range("T3").Value = -4.3292396314701262
$ws.Range("U3").Value = 0.6727781242892189
$ws.Range("V3").Value = -0.42245506744911848
$ws.Range("B4").Value = -24.449606182824091
$ws.Range("C4").Value = -40.894267415818618
$ws.Range("D4").Value = -43.453841318761533
$ws.Range("I4").Value = -9.440410818508397
$ws.Range("J4").Value = -14.77015141898632
$ws.Range("K4").Value = -49.796006858899823
$ws.Range("BF4").Value = 17.58027686899187
$ws.Range("BG4").Value = 27.06795512623556
$ws.Range("BH4").Value = 51.510215141852477
$ws.Range("B5").Value = -7.1650807761036894
$ws.Range("C5").Value = -32.399538106235561
$ws.Range("D5").Value = -42.818060953160163
$ws.Range("BF5").Value = 18.90691327950827
$ws.Range("BG5").Value = 34.965145803913977
$ws.Range("BH5").Value = 32.532308977687258
$ws.Range("B6").Value = -1.8238031628388049
$ws.Range("C6").Value = -11.7519042437432
$ws.Range("D6").Value = 6.736007924715202
$ws.Range("W6").Value = -20.44751796757479
$ws.Range("X6").Value = -26.90802918534623
$ws.Range("Y6").Value = -27.081138040042141
$ws.Range("BB6").Value = -15.2778357976061
$ws.Range("BC6").Value = -20.144404332129959
$ws.Range("BD6").Value = -14.838455459110699
$ws.Range("BE6").Value = -3.3835298433597569
$ws.Range("B7").Value = -17.947313986737239
$ws.Range("C7").Value = -21.9350716014205
$ws.Range("D7").Value = -33.048206993981367
$ws.Range("BF7").Value = -53.055137454773977
$ws.Range("BG7").Value = -53.52959593848631
$ws.Range("BH7").Value = -59.656675044403883
$ws.Range("B8").Value = 26.319883064462051
$ws.Range("C8").Value = 27.57800554957813
$ws.Range("D8").Value = 48.1206179268612
$ws.Range("W8").Value = 9.6642561132591549
$ws.Range("X8").Value = -6.5183587764599134
$ws.Range("Y8").Value = -41.336451970090948
$ws.Range("BF8").Value = -11.53986344966407
$ws.Range("BG8").Value = -32.17778172115927
$ws.Range("BH8").Value = -19.85252410663643
$ws.Range("B9").Value = 19.35075473350701
$ws.Range("C9").Value = 33.087239462374527
$ws.Range("D9").Value = 31.41750463941997
$ws.Range("I9").Value = 20.039018844034459
$ws.Range("J9").Value = 23.9519672932873
$ws.Range("K9").Value = 27.813100130630058
$ws.Range("BF9").Value = -50.511412257806029
$ws.Range("BG9").Value = -61.804027576957729
$ws.Range("BH9").Value = -56.713952165126457
$ws.Range("B10").Value = -7.5087363226424104
$ws.Range("C10").Value = -33.00258929100864
$ws.Range("D10").Value = -50.510298756726669
$ws.Range("BF10").Value = 18.407413366813451
$ws.Range("BG10").Value = 34.02276011560695
$ws.Range("BH10").Value = 39.77531812838739
$ws.Range("B11").Value = -10.31598945770509
$ws.Range("C11").Value = -5.4824511906105373
$ws.Range("D11").Value = 1.249801408674466
$ws.Range("BB11").Value = -26.082407821481169
$ws.Range("BC11").Value = -6.9235902569235872
$ws.Range("BD11").Value = -2.1711872236520851
$ws.Range("BE11").Value = 7.2154502242127938
$ws.Range("B12").Value = 5.0429490633617524
$ws.Range("C12").Value = 5.0529405598070936
$ws.Range("D12").Value = 8.2530539278975041
$ws.Range("W12").Value = -71.843480374286059
$ws.Range("X12").Value = -47.199908847490462
$ws.Range("Y12").Value = -31.849285404937209
$ws.Range("AU12").Value = -15.72026300059771
$ws.Range("AV12").Value = -28.23406734956756
$ws.Range("AW12").Value = -8.8542558657304351
$ws.Range("AX12").Value = -0.52036123140966928
$ws.Range("B13").Value = 227.96812499999999
$ws.Range("C13").Value = 62.112916666666663
$ws.Range("D13").Value = 21.20366666666667
$ws.Range("E13").Value = -318.38474999999988
$ws.Range("F13").Value = -210.94916666666671
$ws.Range("G13").Value = -126.3341666666667
$ws.Range("H13").Value = -34.415335393772907
$ws.Range("I13").Value = 26.952354166666641
$ws.Range("J13").Value = -33.801958333333317
$ws.Range("K13").Value = -76.913083333333333
$ws.Range("L13").Value = -99.588124999999991
$ws.Range("M13").Value = -85.73566666666666
$ws.Range("N13").Value = -71.307791666666688
$ws.Range("O13").Value = -15.030927083333349
$ws.Range("P13").Value = 1.2889791666666379
$ws.Range("Q13").Value = -170.33937499999999
$ws.Range("R13").Value = -247.55083333333329
$ws.Range("S13").Value = -316.52962500000001
$ws.Range("T13").Value = -295.58808333333337
$ws.Range("U13").Value = -169.3099
$ws.Range("V13").Value = -28.445566666666661
$ws.Range("W13").Value = 6.2335208333333352
$ws.Range("X13").Value = -8.5939999999999941
$ws.Range("Y13").Value = -15.28466666666667
$ws.Range("Z13").Value = -24.885000000000002
$ws.Range("AA13").Value = -17.87908333333333
$ws.Range("AB13").Value = -2.0436250000000018
$ws.Range("AC13").Value = 6.5324583333333344
$ws.Range("AD13").Value = -13.73114583333335
$ws.Range("AE13").Value = -123.245125
$ws.Range("AF13").Value = -196.1596666666666
$ws.Range("AG13").Value = -185.44450000000001
$ws.Range("AH13").Value = -137.3410833333333
$ws.Range("AI13").Value = -116.1660833333333
$ws.Range("AJ13").Value = -61.337217803030313
$ws.Range("AK13").Value = -61.819657012195172
$ws.Range("AL13").Value = -142.48065931372551
$ws.Range("AM13").Value = -145.23031666666671
$ws.Range("AN13").Value = -121.94275
$ws.Range("AO13").Value = -183.6809166666666
$ws.Range("AP13").Value = -118.0817361111111
$ws.Range("AQ13").Value = -3.7462749999999998
$ws.Range("AR13").Value = -26.219791666666669
$ws.Range("AS13").Value = -60.007375000000003
$ws.Range("AT13").Value = -82.906583333333344
$ws.Range("AU13").Value = -123.3715
$ws.Range("AV13").Value = -84.723416666666665
$ws.Range("AW13").Value = -35.331500000000013
$ws.Range("AX13").Value = 5.6622083333333437
$ws.Range("AY13").Value = -63.74647916666671
$ws.Range("AZ13").Value = -65.493208333333328
$ws.Range("BA13").Value = -70.894750000000016
$ws.Range("BB13").Value = -98.911633333333384
$ws.Range("BC13").Value = -83.691333333333318
$ws.Range("BD13").Value = -53.405583333333347
$ws.Range("BE13").Value = -43.870229166666682
$ws.Range("BF13").Value = -34.887395833333358
$ws.Range("BG13").Value = -103.70399999999999
$ws.Range("BH13").Value = -145.3764166666667
$ws.Range("BI13").Value = -243.53691666666671
$ws.Range("BJ13").Value = -190.82149999999999
$ws.Range("BK13").Value = -113.7895833333333
$ws.Range("BL13").Value = -70.565645833333335

# --- 3) Clear cells that were removed entirely ---
$ws.Range("I2").ClearContents()
$ws.Range("N2").ClearContents()
$ws.Range("O2").ClearContents()
$ws.Range("W2").ClearContents()
$ws.Range("AA2").ClearContents()
$ws.Range("AB2").ClearContents()
$ws.Range("AC2").ClearContents()
$ws.Range("AR2").ClearContents()
$ws.Range("AS2").ClearContents()
$ws.Range("AT2").ClearContents()
$ws.Range("AU2").ClearContents()
$ws.Range("AV2").ClearContents()
$ws.Range("AW2").ClearContents()
$ws.Range("AX2").ClearContents()
$ws.Range("I6").ClearContents()
$ws.Range("J6").ClearContents()
$ws.Range("K6").ClearContents()
$ws.Range("L6").ClearContents()
$ws.Range("M6").ClearContents()
$ws.Range("N6").ClearContents()
$ws.Range("O6").ClearContents()
$ws.Range("I7").ClearContents()
$ws.Range("J7").ClearContents()
$ws.Range("K7").ClearContents()
$ws.Range("L7").ClearContents()
$ws.Range("M7").ClearContents()
$ws.Range("N7").ClearContents()
$ws.Range("O7").ClearContents()
$ws.Range("W7").ClearContents()
$ws.Range("X7").ClearContents()
$ws.Range("Y7").ClearContents()
$ws.Range("Z7").ClearContents()
$ws.Range("AA7").ClearContents()
$ws.Range("AB7").ClearContents()
$ws.Range("AC7").ClearContents()
$ws.Range("AD7").ClearContents()
$ws.Range("AE7").ClearContents()
$ws.Range("AF7").ClearContents()
$ws.Range("AG7").ClearContents()
$ws.Range("AH7").ClearContents()
$ws.Range("AI7").ClearContents()
$ws.Range("AJ7").ClearContents()
$ws.Range("AK7").ClearContents()
$ws.Range("AL7").ClearContents()
$ws.Range("AM7").ClearContents()
$ws.Range("AN7").ClearContents()
$ws.Range("AO7").ClearContents()
$ws.Range("AP7").ClearContents()
$ws.Range("AQ7").ClearContents()

# --- 4) Update conditional formatting applies-to range (B2:BL13 -> A1:BL13) ---
$fc = $ws.Range("B2:BL13").FormatConditions.Item(1)
$fc.ModifyAppliesToRange($ws.Range("A1:BL13"))

# --- 5) Update sheet view: reset scroll/top-left cell and change selection ---
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.ScrollRow = 1
[void]$ws.Range("BF3:BL3").Select()
